# Updates the cryptocurrency price/volume snapshot (Price = column D,
# Volume(1h) = column E) for rows 2-47 of Sheet1, matching the refreshed
# scrape values from the "Updated symbol list" GitHub Actions commit.
#
# Each target cell stores its value as literal text (e.g. "304.71",
# "2.56%") in the source workbook, rather than as a genuine Excel number
# or percentage. Assigning a numeric-looking string straight to .Value
# would let Excel auto-convert it to a real number (and stamp a percent
# number format on the "%" cells), so each cell is briefly switched to
# the Text ("@") number format before the value is written, then
# restored to the workbook's default "Normal" style so no stray
# formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($addr, $value)
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2
Set-TextValue "D2" "304.71"
Set-TextValue "E2" "2.56%"

# Row 3
Set-TextValue "D3" "31.96"
Set-TextValue "E3" "0.74%"

# Row 4
Set-TextValue "D4" "5.187"
Set-TextValue "E4" "0.70%"

# Row 5
Set-TextValue "D5" "0.07469"
Set-TextValue "E5" "-0.03%"

# Row 6
Set-TextValue "D6" "2.408"
Set-TextValue "E6" "46.05%"

# Row 7
Set-TextValue "D7" "7.998"
Set-TextValue "E7" "2.94%"

# Row 8
Set-TextValue "D8" "3.862"
Set-TextValue "E8" "1.36%"

# Row 9
Set-TextValue "D9" "0.9187"
Set-TextValue "E9" "-0.62%"

# Row 10
Set-TextValue "D10" "0.1736"
Set-TextValue "E10" "1.54%"

# Row 11
Set-TextValue "D11" "0.07666"
Set-TextValue "E11" "1.13%"

# Row 12
Set-TextValue "D12" "0.08208"
Set-TextValue "E12" "3.75%"

# Row 13
Set-TextValue "D13" "0.03009"
Set-TextValue "E13" "0.39%"

# Row 14
Set-TextValue "D14" "0.09931"
Set-TextValue "E14" "0.49%"

# Row 15
Set-TextValue "D15" "0.001528"
Set-TextValue "E15" "2.34%"

# Row 16
Set-TextValue "D16" "0.006134"
Set-TextValue "E16" "-1.49%"

# Row 17
Set-TextValue "D17" "3.501"
Set-TextValue "E17" "1.70%"

# Row 18
Set-TextValue "E18" "-0.06%"

# Row 19
Set-TextValue "D19" "0.3262"
Set-TextValue "E19" "-0.93%"

# Row 20
Set-TextValue "D20" "0.1337"
Set-TextValue "E20" "-0.83%"

# Row 21
Set-TextValue "D21" "4.653"
Set-TextValue "E21" "1.79%"

# Row 22
Set-TextValue "D22" "0.04601"
Set-TextValue "E22" "-1.30%"

# Row 23
Set-TextValue "E23" "0.72%"

# Row 24
Set-TextValue "E24" "3.10%"

# Row 25
Set-TextValue "D25" "0.004526"
Set-TextValue "E25" "2.54%"

# Row 26
Set-TextValue "D26" "0.0001298"
Set-TextValue "E26" "-7.44%"

# Row 27
Set-TextValue "D27" "0.0002738"
Set-TextValue "E27" "51.42%"

# Row 39
Set-TextValue "D39" "0.01781"
Set-TextValue "E39" "7.83%"

# Row 40
Set-TextValue "D40" "0.04559"
Set-TextValue "E40" "0.88%"

# Row 41
Set-TextValue "D41" "0.007374"
Set-TextValue "E41" "6.11%"

# Row 42
Set-TextValue "D42" "0.1363"
Set-TextValue "E42" "1.56%"

# Row 43
Set-TextValue "D43" "0.002177"
Set-TextValue "E43" "5.49%"

# Row 44
Set-TextValue "D44" "0.01081"
Set-TextValue "E44" "-19.90%"

# Row 45
Set-TextValue "D45" "0.00006463"
Set-TextValue "E45" "6.14%"

# Row 46
Set-TextValue "E46" "15.26%"

# Row 47
Set-TextValue "D47" "0.009883"
Set-TextValue "E47" "-19.39%"
